# Fix stepper import tax type
# Adds a new "tax_type" column (E) to the "Productos" sheet, with a mix of
# percentage values (for GRAVADO rates) and free-text labels (Exento /
# Exonerado / GRAVADO_15 / GRAVADO_18), matching the bulk-import template.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productos")

# Copy the existing header/text style (as used by columns A-D) onto the new
# column so the new cells reuse the same font/alignment instead of creating
# brand-new style entries.
$ws.Range("A1").Copy()
$ws.Range("E1:E10").PasteSpecial(-4122)

# Header
$ws.Range("E1").Value = "tax_type"

# Percentage (numeric) tax rates -> formatted as 0%
$percentCells = @("E2", "E3", "E4", "E7", "E10")
foreach ($cell in $percentCells) {
    $ws.Range($cell).NumberFormat = "0%"
}

$ws.Range("E2").Value = 0.0
$ws.Range("E3").Value = 0.18
$ws.Range("E4").Value = 0.15
$ws.Range("E7").Value = 0.0
$ws.Range("E10").Value = 0.18

# Free-text tax type labels
$ws.Range("E5").Value = "Exento"
$ws.Range("E6").Value = "Exonerado"
$ws.Range("E8").Value = "GRAVADO_15"
$ws.Range("E9").Value = "GRAVADO_18"
